$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3, shifting existing row 3 (and below) down to row 4.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted A3 with the same URL as A2.
$ws.Range("A3").Value = "https://www.imdb.com/title/tt13622776/?ref_=hm_top_tt_i_1"
